# AudaciaBall Planning.xlsx update
# "publish API, Add vue components, add CORS to API"
#
# 1) Planning sheet: mark a couple of existing tasks as Done + fill in their
#    time spent, insert two new "Deploy" tasks for Azure (DB + API), and
#    replace the old "Simply deploiement with Docker" task with
#    "Deploy UI on web".
# 2) Used sheet: add a new "Cors" dependency row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Planning"
# ---------------------------------------------------------------------
$planning = $wb.Worksheets.Item("Planning")

# "Add UI components" (row 12) is now finished -> mark Done + log 2h.
$planning.Range("D12").Value = "Yes"
$planning.Range("E12").Value = 2

# "Test API" (row 15) is now finished -> mark Done + log 1h.
$planning.Range("D15").Value = "Yes"
$planning.Range("E15").Value = 1

# Make room for two new "Deploy" rows right after "Test API" (row 15),
# pushing every following row down by two (formulas/refs auto-adjust).
$planning.Rows.Item(16).Insert()
$planning.Rows.Item(16).Insert()

# Copy the formatting (borders etc.) of the row just above down onto the
# two freshly inserted, still-blank rows.
$planning.Range("A15:E15").Copy()
$planning.Range("A16:E17").PasteSpecial(-4122)
$planning.Application.CutCopyMode = $false

$planning.Range("A16").Value = 10.1
$planning.Range("B16").Value = "Deploy"
$planning.Range("C16").Value = "Create DB on Azure"
$planning.Range("D16").Value = "Yes"
$planning.Range("E16").Value = 0.5

$planning.Range("A17").Value = 10.2
$planning.Range("B17").Value = "Deploy"
$planning.Range("C17").Value = "Deploy API on Azure"
$planning.Range("D17").Value = "Yes"
$planning.Range("E17").Value = 2

# Used sheet gets the CORS dependency row now (so new shared strings come
# out in the same order as the original edit).
$used = $wb.Worksheets.Item("Used")
$used.Range("A5").Value = "Cors"
$used.Range("B5").Value = "installed with install-package"

# The old Docker deployment task became a "Deploy UI on web" task (same
# row, same "Deploy" category, after the two rows inserted above it now
# this used to be row 18 and is row 20).
$planning.Range("C20").Value = "Deploy UI on web"

# Keep the selection/scroll position close to where the author left off,
# making sure "Planning" ends up as the active sheet again.
$used.Range("B6").Select()
$planning.Activate()
$planning.Range("C17").Select()
